$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Layout")

# Row 8: turn on P8:V8 (0 -> 1)
$ws.Range("P8:V8").Value = 1

# Row 13: turn off several blocks (1 -> 0)
$ws.Range("C13:F13").Value = 0
$ws.Range("J13:K13").Value = 0
$ws.Range("R13:T13").Value = 0
$ws.Range("AA13:AB13").Value = 0
$ws.Range("AF13:AI13").Value = 0

# Row 18: turn off two blocks (1 -> 0)
$ws.Range("C18:F18").Value = 0
$ws.Range("AF18:AI18").Value = 0

# Row 23: mix of toggles
$ws.Range("C23:E23").Value = 0
$ws.Range("G23:I23").Value = 1
$ws.Range("M23:O23").Value = 0
$ws.Range("W23:Y23").Value = 0
$ws.Range("AC23:AE23").Value = 1
$ws.Range("AG23:AI23").Value = 0

# Update the active cell selection to match the editor's last position
$ws.Range("X20").Select()
